$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.060.64"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.376.95"
$ws.Range("E3").Value = "  +6.64%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.43%  "
$ws.Range("E7").Value = "  +2.63%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.663"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.14%  "
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.45%  "
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "2.736.58"
$ws.Range("E16").Value = "  +6.68%  "
$ws.Range("D17").Value = "2.400.62"
$ws.Range("E17").Value = "  +7.87%  "
$ws.Range("D18").Value = "43.028.47"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.02%  "
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "276.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.11%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.20%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0918"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("E34").Value = "  +2.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.52%  "
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0365"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +18.07%  "
$ws.Range("E41").Value = "  +21.87%  "
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +20.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +64.70%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("E48").Value = "  +11.21%  "
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").Value = "1.608.99"
$ws.Range("E51").Value = "  +12.44%  "
